$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("London_Eq")

# Row 3: TRD-DEMO1 (matched/affirmed)
$row3 = @("TRD-DEMO1","LEG-DEMO","ALLOC-DEMO1","DEMO-PORT","ISIN-DEMO1","USD","25/08/2024","27/08/2024",0.0092,100000,99500,0.01,175,1500,166.4,"READY","AFFIRMED","DemoCounterparty","DemoBroker","C-JPY","GLOBAL","GLOBAL-DEMO","emea_positions.xlsx","MATCH-CONFIRMED",0.96)

# Row 4: TRD-DEMO2 (synthetic mismatch)
$row4 = @("TRD-DEMO2","LEG-DEMO","ALLOC-DEMO2","DEMO-PORT","ISIN-DEMO2","JPY","21/08/2024","23/08/2024",0.0092,255500,255000,0.01,175,1500,170.3,"READY","MATCHED","EuroPrime","BrokerLondon","C-JPY","EMEA","EMEA-DEMO","emea_positions.xlsx","Synthetic mismatch",0.96)

for ($col = 1; $col -le $row3.Length; $col++) {
    $ws.Cells.Item(3, $col).Value = $row3[$col - 1]
}

for ($col = 1; $col -le $row4.Length; $col++) {
    $ws.Cells.Item(4, $col).Value = $row4[$col - 1]
}
